$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the text in B73: add "prava, " and ", obrazku"
$ws.Range("B73").Value = "Psaní - formulare, otazky, role, prava, ...; reseni tabulek, obrazku"

# Update C73 value from 3 to 6
$ws.Range("C73").Value = 6

# Update the selection on the active sheet to C74
$ws.Range("C74").Select()
